$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.230.35"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3
$ws.Range("D3").Value = "1.906.71"
$ws.Range("E3").Value = "  +0.23%  "

# Row 4
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.72"
$ws.Range("E5").Value = "  +0.60%  "

# Row 6
$ws.Range("E6").Value = "  +0.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5254"
$ws.Range("E7").Value = "  +0.47%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3822"
$ws.Range("E8").Value = "  +1.66%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07306"
$ws.Range("E9").Value = "  +0.83%  "

# Row 10
$ws.Range("E10").Value = "  +2.18%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9058"
$ws.Range("E11").Value = "  +0.43%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08089"
$ws.Range("E12").Value = "  -4.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.00"
$ws.Range("E13").Value = "  +1.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.369"
$ws.Range("E14").Value = "  +1.55%  "

# Row 15
$ws.Range("D15").Value = "1.772.16"
$ws.Range("E15").Value = "  -6.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.14%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008674"
$ws.Range("E17").Value = "  +0.57%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.75"
$ws.Range("E18").Value = "  +1.61%  "

# Row 20
$ws.Range("D20").Value = "27.271.57"
$ws.Range("E20").Value = "  +0.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.121"
$ws.Range("E21").Value = "  +1.11%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.82"
$ws.Range("E22").Value = "  +2.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.486"
$ws.Range("E23").Value = "  +0.95%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.345"
$ws.Range("E24").Value = "  +2.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.15"
$ws.Range("E25").Value = "  +1.91%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.741"
$ws.Range("E27").Value = "  -0.66%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.79"
$ws.Range("E28").Value = "  +1.69%  "

# Row 29
$ws.Range("E29").Value = "  +0.93%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.881"
$ws.Range("E30").Value = "  -0.12%  "

# Row 31
$ws.Range("E31").Value = "  -0.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8257"
$ws.Range("E32").Value = "  +2.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05081"
$ws.Range("E33").Value = "  +0.57%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.232"
$ws.Range("E34").Value = "  -0.25%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.994"
$ws.Range("E35").Value = "  +1.70%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.365"
$ws.Range("E36").Value = "  -2.35%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.733"
$ws.Range("E37").Value = "  +4.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5753"
$ws.Range("E38").Value = "  +0.68%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02004"
$ws.Range("E39").Value = "  +0.52%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.084"
$ws.Range("E40").Value = "  +0.84%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.056"
$ws.Range("E41").Value = "  +0.42%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.619"
$ws.Range("E42").Value = "  -0.13%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.15"
$ws.Range("E43").Value = "  +0.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1524"
$ws.Range("E44").Value = "  +0.64%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4935"
$ws.Range("E45").Value = "  +1.63%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.20"
$ws.Range("E46").Value = "  +0.90%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.18%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.642"
$ws.Range("E48").Value = "  +1.56%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.66"
$ws.Range("E49").Value = "  +3.16%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.20"
$ws.Range("E50").Value = "  +0.28%  "

# Row 51
$ws.Range("E51").Value = "  +0.41%  "
